# Add a new "MVP done" column (inserted before the existing "Change detection"
# column, i.e. before column E) and populate it for the sources whose MCP
# integration is already working, plus fix a couple of priority values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; this shifts the old E:K (Change detection .. Latency)
# right to F:L, carrying along their formatting/column widths.
$ws.Range("E1").EntireColumn.Insert()

# New column header
$ws.Range("E1").Value = "MVP done"

# Mark rows whose MCP integration already has a working proof of concept
$ws.Range("E2").Value = "Yes"
$ws.Range("E3").Value = "Yes"
$ws.Range("E4").Value = "Yes"
$ws.Range("E8").Value = "Write only!"
$ws.Range("E10").Value = "Yes"
$ws.Range("E15").Value = "Yes"
$ws.Range("E23").Value = "Yes"
$ws.Range("E24").Value = "Yes"
$ws.Range("E35").Value = "Yes"

# Priority corrections
$ws.Range("B5").Value = 10
$ws.Range("B6").Value = 6

# Restore the selection to where the user last left off
$ws.Range("E39").Select()
